$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dAddr = "D$r"
    $eAddr = "E$r"
    $fAddr = "F$r"
    $gAddr = "G$r"

    # Target mapping (value flow): D_new = G_old ; E_new = F_old ; F_new = D_old ; G_new = E_old
    # Rotate using one helper cell: temp = G; G = E; E = F; F = D; D = temp
    $ws.Range($gAddr).Copy($ws.Range("Z1"))
    $ws.Range($eAddr).Copy($ws.Range($gAddr))
    $ws.Range($fAddr).Copy($ws.Range($eAddr))
    $ws.Range($dAddr).Copy($ws.Range($fAddr))
    $ws.Range("Z1").Copy($ws.Range($dAddr))
}

$ws.Range("Z1").ClearContents()
Write-Host "done"
